$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 181 - this shifts all the
# existing rows 181..282 down to 182..283 and grows the used range to
# A1:R283 (matching the new <dimension> in the diff).
$ws.Rows.Item(181).Insert()

# Populate the freshly inserted row 181 with the new price record.
$ws.Range("A181").Value = 4
$ws.Range("B181").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C181").Value = "Los Lagos"
$ws.Range("D181").Value = 44529
$ws.Range("E181").Value = 10
$ws.Range("F181").Value = 100114001
$ws.Range("G181").Value = "Papa"
$ws.Range("H181").Value = "Pehuenche"
$ws.Range("I181").Value = "1a nueva(o)"
$ws.Range("J181").Value = 300
$ws.Range("K181").Value = 13000
$ws.Range("L181").Value = 14000
$ws.Range("M181").Value = 13500
$ws.Range("N181").Value = "$/saco 25 kilos"
$ws.Range("O181").Value = "Región de La Araucanía"
$ws.Range("P181").Value = 540
$ws.Range("Q181").Value = 25
$ws.Range("R181").Value = "Hortaliza"
